$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 82
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 11
$ws.Range("E52").Value = 24
$ws.Range("F52").Value = 94
$ws.Range("G52").Value = 118
